$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.831.64"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "2.649.58"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'600.86"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'155.95"
$ws.Range("E6").Value = "  +4.02%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "2.646.71"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").Value = "'0.140"
$ws.Range("E10").Value = "  +13.24%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("E15").Value = "  +6.07%  "
$ws.Range("D16").Value = "3.132.33"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "68.692.91"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").Value = "2.641.49"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").Value = "'365.66"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'4.91"
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("E24").Value = "  +4.51%  "
$ws.Range("E25").Value = "  +10.20%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'10.09"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("E28").Value = "  +7.74%  "
$ws.Range("D30").Value = "'583.57"
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  +4.52%  "
$ws.Range("E33").Value = "  +4.83%  "
$ws.Range("E34").Value = "  +3.08%  "
$ws.Range("D35").Value = "'0.132"
$ws.Range("E35").Value = "  +5.53%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +4.17%  "
$ws.Range("D38").Value = "'160.44"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("E39").Value = "  +4.42%  "
$ws.Range("D40").Value = "'19.35"
$ws.Range("E40").Value = "  +1.82%  "
$ws.Range("E41").Value = "  +0.85%  "
$ws.Range("E42").Value = "  +3.46%  "
$ws.Range("E43").Value = "  +7.08%  "
$ws.Range("D44").Value = "'17.72"
$ws.Range("E44").Value = "  +5.45%  "
$ws.Range("E45").Value = "  +13.01%  "
$ws.Range("D46").Value = "'40.77"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "'156.28"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("E50").Value = "  +3.30%  "
$ws.Range("E51").Value = "  +1.44%  "
